# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Golem_Profits workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 448.75
$ws.Range("I52").Value = 448.75
$ws.Range("K52").Value = 1346.25
$ws.Range("M52").Value = -1186.25
$ws.Range("H64").Value = 4249.5
$ws.Range("J64").Value = 4249.5
$ws.Range("L64").Value = 4249.5
$ws.Range("N64").Value = -4745.5
$ws.Range("H67").Value = 4249.5
$ws.Range("J67").Value = 4249.5
$ws.Range("L67").Value = 4249.5
$ws.Range("N67").Value = -5965.5
$ws.Range("H98").Value = 1265.8334
$ws.Range("I98").Value = 1265.8334
$ws.Range("K98").Value = 1265.8334
$ws.Range("M98").Value = 232.1666
$ws.Range("H122").Value = 1265.8334
$ws.Range("I122").Value = 1265.8334
$ws.Range("K122").Value = 3797.5002
$ws.Range("M122").Value = -1347.5002
$ws.Range("H129").Value = 1865.5
$ws.Range("J129").Value = 1999
$ws.Range("L129").Value = 5997
$ws.Range("N129").Value = -15997
$ws.Range("H132").Value = 934.8333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 79.333336
$ws.Range("I2").Value = 72.5
$ws.Range("K2").Value = 72.5
$ws.Range("M2").Value = 40.5
$ws.Range("H4").Value = 263.76923
$ws.Range("I4").Value = 305.6
$ws.Range("K4").Value = 305.6
$ws.Range("M4").Value = -189.6
$ws.Range("H32").Value = 2933.2307
$ws.Range("I32").Value = 2933.2307
$ws.Range("K32").Value = 2933.2307
$ws.Range("M32").Value = -2646.2307
$ws.Range("H116").Value = 79.333336
$ws.Range("I116").Value = 72.5
$ws.Range("K116").Value = 72.5
$ws.Range("M116").Value = 2221.5
$ws.Range("H132").Value = 2422
$ws.Range("I132").Value = 2422
$ws.Range("K132").Value = 7266
$ws.Range("M132").Value = -4736

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 79.333336
$ws.Range("I3").Value = 72.5
$ws.Range("K3").Value = 72.5
$ws.Range("M3").Value = 41.5
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H105").Value = 2009
$ws.Range("I105").Value = 2008.5
$ws.Range("K105").Value = 2008.5
$ws.Range("M105").Value = -261.5
$ws.Range("H107").Value = 69966.664
$ws.Range("I107").Value = 83160
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 83160
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -81240
$ws.Range("N107").Value = -7840
$ws.Range("H134").Value = 1415.6
$ws.Range("I134").Value = 898.25
$ws.Range("K134").Value = 2694.75
$ws.Range("M134").Value = -159.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 612.4
$ws.Range("J22").Value = 417.14285
$ws.Range("L22").Value = 417.14285
$ws.Range("N22").Value = -1117.14285
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H56").Value = 72500
$ws.Range("I56").Value = 72500
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 72500
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -71655
$ws.Range("N56").ClearContents()
$ws.Range("H82").Value = 56000
$ws.Range("J82").Value = 56000
$ws.Range("L82").Value = 56000
$ws.Range("N82").Value = -56722
$ws.Range("H85").Value = 56000
$ws.Range("J85").Value = 56000
$ws.Range("L85").Value = 56000
$ws.Range("N85").Value = -58496
$ws.Range("H99").Value = 1669000
$ws.Range("I99").Value = 1669000
$ws.Range("K99").Value = 1669000
$ws.Range("M99").Value = -1667502
$ws.Range("H105").Value = 504.375
$ws.Range("I105").Value = 554.8570999999999
$ws.Range("J105").Value = 151
$ws.Range("K105").Value = 554.8570999999999
$ws.Range("L105").Value = 151
$ws.Range("M105").Value = 1192.1429
$ws.Range("N105").Value = -3645
$ws.Range("H126").Value = 1669000
$ws.Range("I126").Value = 1669000
$ws.Range("K126").Value = 5007000
$ws.Range("M126").Value = -5004530

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1868.909
$ws.Range("I68").Value = 1686.75
$ws.Range("J68").Value = 1973
$ws.Range("K68").Value = 5060.25
$ws.Range("L68").Value = 5919
$ws.Range("M68").Value = -4249.25
$ws.Range("N68").Value = -7541
$ws.Range("H71").Value = 1868.909
$ws.Range("I71").Value = 1686.75
$ws.Range("J71").Value = 1973
$ws.Range("K71").Value = 15180.75
$ws.Range("L71").Value = 17757
$ws.Range("M71").Value = -11124.75
$ws.Range("N71").Value = -25869
$ws.Range("H92").Value = 410.5
$ws.Range("I92").Value = 410.5
$ws.Range("K92").Value = 1231.5
$ws.Range("M92").Value = 16.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 11999
$ws.Range("J34").Value = 11999
$ws.Range("L34").Value = 11999
$ws.Range("N34").Value = -12535
$ws.Range("H76").Value = 11999
$ws.Range("J76").Value = 11999
$ws.Range("L76").Value = 11999
$ws.Range("N76").Value = -12629
$ws.Range("H79").Value = 11999
$ws.Range("J79").Value = 11999
$ws.Range("L79").Value = 11999
$ws.Range("N79").Value = -14183
$ws.Range("H80").Value = 5491
$ws.Range("I80").Value = 5499
$ws.Range("J80").Value = 5487
$ws.Range("K80").Value = 5499
$ws.Range("L80").Value = 5487
$ws.Range("M80").Value = -4501
$ws.Range("N80").Value = -7483
$ws.Range("H83").Value = 5491
$ws.Range("I83").Value = 5499
$ws.Range("J83").Value = 5487
$ws.Range("K83").Value = 27495
$ws.Range("L83").Value = 27435
$ws.Range("M83").Value = -22503
$ws.Range("N83").Value = -37419
$ws.Range("H102").Value = 1506
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 5104.8237
$ws.Range("I122").Value = 3779.5
$ws.Range("J122").Value = 6998.143
$ws.Range("K122").Value = 11338.5
$ws.Range("L122").Value = 20994.429
$ws.Range("M122").Value = -8888.5
$ws.Range("N122").Value = -25894.429
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820
$ws.Range("H126").Value = 3120.8333
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1648.75
$ws.Range("I132").Value = 1511.5
$ws.Range("J132").Value = 2609.5
$ws.Range("K132").Value = 4534.5
$ws.Range("L132").Value = 7828.5
$ws.Range("M132").Value = -2004.5
$ws.Range("N132").Value = -12888.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3363.75
$ws.Range("I68").Value = 3262.3333
$ws.Range("J68").Value = 3465.1667
$ws.Range("K68").Value = 3262.3333
$ws.Range("L68").Value = 3465.1667
$ws.Range("M68").Value = -2513.3333
$ws.Range("N68").Value = -4963.1667
$ws.Range("H71").Value = 3363.75
$ws.Range("I71").Value = 3262.3333
$ws.Range("J71").Value = 3465.1667
$ws.Range("K71").Value = 16311.6665
$ws.Range("L71").Value = 17325.8335
$ws.Range("M71").Value = -12567.6665
$ws.Range("N71").Value = -24813.8335

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 213579.2
$ws.Range("I2").Value = 343299.66
$ws.Range("J2").Value = 18998.5
$ws.Range("K2").Value = 343299.66
$ws.Range("L2").Value = 18998.5
$ws.Range("M2").Value = -343187.66
$ws.Range("N2").Value = -19222.5
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H126").Value = 5212.125
$ws.Range("I126").Value = 4020.6
$ws.Range("K126").Value = 12061.8
$ws.Range("M126").Value = -9591.799999999999
